$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, avoiding auto-numeric conversion of values like "185.45"
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.395.10'
$ws.Range("E2").Value = '  -4.71%  '

$ws.Range("D3").Value = '3.386.95'
$ws.Range("E3").Value = '  -6.71%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '185.45'
$ws.Range("E5").Value = '  -8.70%  '

$ws.Range("D6").Value = '531.26'
$ws.Range("E6").Value = '  -7.82%  '

$ws.Range("D7").Value = '0.602'
$ws.Range("E7").Value = '  -3.66%  '

$ws.Range("D8").Value = '3.387.36'
$ws.Range("E8").Value = '  -6.56%  '

$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").Value = '0.627'
$ws.Range("E10").Value = '  -9.34%  '

$ws.Range("D11").Value = '57.73'
$ws.Range("E11").Value = '  -6.00%  '

$ws.Range("D12").Value = '0.133'
$ws.Range("E12").Value = '  -12.13%  '

$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -11.94%  '

$ws.Range("D14").Value = '9.30'
$ws.Range("E14").Value = '  -8.32%  '

$ws.Range("D15").Value = '3.911.09'
$ws.Range("E15").Value = '  -6.98%  '

$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("D17").Value = '3.378.07'
$ws.Range("E17").Value = '  -7.06%  '

$ws.Range("D18").Value = '65.122.41'
$ws.Range("E18").Value = '  -4.89%  '

$ws.Range("D19").Value = '17.42'
$ws.Range("E19").Value = '  -9.02%  '

$ws.Range("D20").Value = '11.08'
$ws.Range("E20").Value = '  -11.61%  '

$ws.Range("D21").Value = '0.970'
$ws.Range("E21").Value = '  -10.26%  '

$ws.Range("D22").Value = '373.96'
$ws.Range("E22").Value = '  -8.18%  '

$ws.Range("D23").Value = '81.29'
$ws.Range("E23").Value = '  -5.68%  '

$ws.Range("D24").Value = '3.73'
$ws.Range("E24").Value = '  -12.20%  '

$ws.Range("D25").Value = '10.83'
$ws.Range("E25").Value = '  -15.89%  '

$ws.Range("D26").Value = '3.80'
$ws.Range("E26").Value = '  -5.08%  '

$ws.Range("D27").Value = '5.84'
$ws.Range("E27").Value = '  -5.31%  '

$ws.Range("E28").Value = '  -10.52%  '

$ws.Range("D29").Value = '11.47'
$ws.Range("E29").Value = '  -9.73%  '

$ws.Range("D30").Value = '8.52'
$ws.Range("E30").Value = '  -9.94%  '

$ws.Range("D31").Value = '29.59'
$ws.Range("E31").Value = '  -7.33%  '

$ws.Range("D32").Value = '667.64'
$ws.Range("E32").Value = '  -1.57%  '

$ws.Range("D33").Value = '6.76'
$ws.Range("E33").Value = '  -14.19%  '

$ws.Range("D34").Value = '11.18'
$ws.Range("E34").Value = '  -9.32%  '

$ws.Range("D35").Value = '61.01'
$ws.Range("E35").Value = '  -4.29%  '

$ws.Range("D36").Value = '0.105'
$ws.Range("E36").Value = '  -9.27%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '36.58'
$ws.Range("E38").Value = '  -13.20%  '

$ws.Range("D39").Value = '0.381'
$ws.Range("E39").Value = '  -10.09%  '

$ws.Range("D40").Value = '0.995'
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("D41").Value = '0.128'
$ws.Range("E41").Value = '  -5.97%  '

$ws.Range("D42").Value = '2.829.20'
$ws.Range("E42").Value = '  -12.12%  '

$ws.Range("D43").Value = '2.77'
$ws.Range("E43").Value = '  -14.48%  '

$ws.Range("D44").Value = '0.0₃0628'
$ws.Range("E44").Value = '  -20.23%  '

$ws.Range("D45").Value = '0.0392'
$ws.Range("E45").Value = '  -6.90%  '

$ws.Range("D46").Value = '2.62'
$ws.Range("E46").Value = '  -9.95%  '

$ws.Range("D47").Value = '2.36'
$ws.Range("E47").Value = '  -13.54%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.125'
$ws.Range("E48").Value = '  -5.78%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '136.66'
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("D50").Value = '2.86'
$ws.Range("E50").Value = '  -7.51%  '

$ws.Range("D51").Value = '2.60'
$ws.Range("E51").Value = '  -5.84%  '
